$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -7
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = -4
$ws.Range("F5").Value = 0
